$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bValues = @(142,142,172,145,142,142,140,143,140,167,167,140,140,140,151,140,140,140,143,140,147,140,140,140,140,140,147,140,151,140,166,166,140,143,167,166,166,170,140,167,167,150,140,140,140,140,161,140,140,140,140,148,140,140,140,140,140,148,140,167,140,151,167,140,150,140,140,140,150,167,140,144,140,140,166,140,140,140,140,144,150,150,150,167,140,140,140,144,140,143,140,143,161,140,139,140,140,148,167,140)
$cValues = @(40,40,80,79,40,40,40,79,40,80,80,40,40,40,80,40,40,40,79,40,79,40,40,40,40,40,79,40,80,40,79,80,40,79,80,79,79,80,40,80,80,79,40,40,40,40,79,40,40,40,40,80,40,40,40,40,40,80,40,80,40,80,80,40,79,40,40,40,79,80,40,80,40,40,79,40,40,40,40,80,79,79,79,80,40,40,40,80,40,79,40,79,79,40,39,40,40,80,80,40)

for ($i = 0; $i -lt $bValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}
